# Rename the existing sheet and add a new one after it, matching the
# target workbook layout: devTestLogin (sheetId 1), LoginTest (sheetId 2).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "devTestLogin"

# Copy (rather than Add) the sheet so the new tab inherits the same
# sheetFormatPr / namespace declarations as the original, matching the
# target worksheet XML more closely.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "LoginTest"
$ws2.Cells.Clear()

# Headers first on both sheets (matches shared-string build order of the
# authored workbook), then the data rows.
$ws1.Cells.Item(1, 1).Value = "USERNAME"
$ws1.Cells.Item(1, 2).Value = "PASSWORD"

$ws2.Cells.Item(1, 1).Value = "username"
$ws2.Cells.Item(1, 2).Value = "password"

# ---- devTestLogin (sheet1) data rows ----
$devData = @(
    @("admin1", "useradmin"),
    @("admin", "useadmin1"),
    @("admina", "useradmina"),
    @("admin", "useradmin"),
    @("admin", "useradmin"),
    @("admin", "useradmin"),
    @("admin", "useradmin"),
    @("admin", "useradmin"),
    @("admin", "useradmin")
)

for ($i = 0; $i -lt $devData.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $devData[$i][0]
    $ws1.Cells.Item($row, 2).Value = $devData[$i][1]
}

$ws1.Range("B9").Select()

# ---- LoginTest (sheet2) data rows ----
$loginData = @(
    @("admin1", "useradmin"),
    @("admin", "useadmin1"),
    @("admina", "useradmina"),
    @("admin", "useradmin")
)

for ($i = 0; $i -lt $loginData.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $loginData[$i][0]
    $ws2.Cells.Item($row, 2).Value = $loginData[$i][1]
}

$ws2.Range("D4").Select()

# Re-select the devTestLogin tab / cell so it's the active sheet+view,
# matching the tabSelected flag in the target sheet1.xml.
$ws1.Activate()
$ws1.Range("B9").Select()
